$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the column-header shared strings: "_old" -> "_FV2310", "_new" -> "_FV2404"
$headers = @(
    "Segmentname_FV2310", "Segmentgruppe_FV2310", "Segment_FV2310", "Datenelement_FV2310",
    "Segment ID_FV2310", "Code_FV2310", "Qualifier_FV2310", "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310", "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404",
    "Segment ID_FV2404", "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404", "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# 2) Turn the used range into an Excel Table ("Table1") with the new headers
$rng = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split after row 1, top-left of the scrolling pane is A2)
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
